# "Generate Report for Archive"
# Refresh the localization status report:
#   - 6423883b... and f2ffb738... move from "Ready for handoff" to "In Translation"
#   - rows for f2ffb738 / b8dc45cb / cd865f8c re-sort so f2ffb738 moves up
#     ahead of b8dc45cb (still behind 6423883b, still ahead of cd865f8c)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name / zh-cn / de-de / Latest Handoff Date
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

$ws.Range("A5").Value = "f2ffb738-5e69-400b-b1f4-4913a1c05516.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "In Translation"

$ws.Range("A6").Value = "b8dc45cb-ddaa-4348-954f-46f3cfc9aa28.md"

$ws.Range("A7").Value = "cd865f8c-ef21-4a4c-a79f-fd30c0e95a66.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": detailed per-file status
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C4").Value = "In Translation"

$ws.Range("A5").Value = "f2ffb738-5e69-400b-b1f4-4913a1c05516.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("D5").Value = "f2ffb738-5e69-400b-b1f4-4913a1c05516.07385779dc10aa2adc103d308e62048088991871.zh-cn.xlf"

$ws.Range("A6").Value = "b8dc45cb-ddaa-4348-954f-46f3cfc9aa28.md"
$ws.Range("D6").Value = "b8dc45cb-ddaa-4348-954f-46f3cfc9aa28.6b59cec54b0baf97f56621a91a457ee5064bd3e3.zh-cn.xlf"

$ws.Range("A7").Value = "cd865f8c-ef21-4a4c-a79f-fd30c0e95a66.md"
$ws.Range("D7").Value = "cd865f8c-ef21-4a4c-a79f-fd30c0e95a66.40988aa5b4a284abcf75269c2d06e1613d9ae202.zh-cn.xlf"

# ---------------------------------------------------------------------------
# Sheet "de-de": detailed per-file status
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C4").Value = "In Translation"

$ws.Range("A5").Value = "f2ffb738-5e69-400b-b1f4-4913a1c05516.md"
$ws.Range("C5").Value = "In Translation"
$ws.Range("D5").Value = "f2ffb738-5e69-400b-b1f4-4913a1c05516.07385779dc10aa2adc103d308e62048088991871.de-de.xlf"

$ws.Range("A6").Value = "b8dc45cb-ddaa-4348-954f-46f3cfc9aa28.md"
$ws.Range("D6").Value = "b8dc45cb-ddaa-4348-954f-46f3cfc9aa28.6b59cec54b0baf97f56621a91a457ee5064bd3e3.de-de.xlf"

$ws.Range("A7").Value = "cd865f8c-ef21-4a4c-a79f-fd30c0e95a66.md"
$ws.Range("D7").Value = "cd865f8c-ef21-4a4c-a79f-fd30c0e95a66.40988aa5b4a284abcf75269c2d06e1613d9ae202.de-de.xlf"
